$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the columns we touch so numeric-looking strings
# (e.g. "230.93", "61.01") are preserved as text, matching the source data
# which is stored as inline strings, not numbers.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "41.444.65"
$ws.Range("E2").Value = "  +4.32%  "
$ws.Range("D3").Value = "2.219.62"
$ws.Range("E3").Value = "  +2.67%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "230.93"
$ws.Range("E5").Value = "  +1.81%  "
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("D7").Value = "61.01"
$ws.Range("E7").Value = "  -2.91%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +2.89%  "
$ws.Range("D10").Value = "58.63"
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("D11").Value = "0.0887"
$ws.Range("E11").Value = "  +5.41%  "
$ws.Range("D12").Value = "0.104"
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("D13").Value = "2.548.97"
$ws.Range("E13").Value = "  +2.65%  "
$ws.Range("D14").Value = "15.65"
$ws.Range("E14").Value = "  -1.42%  "
$ws.Range("D15").Value = "21.83"
$ws.Range("E15").Value = "  +0.29%  "
$ws.Range("D16").Value = "0.799"
$ws.Range("E16").Value = "  -1.05%  "
$ws.Range("D17").Value = "5.55"
$ws.Range("E17").Value = "  +1.12%  "
$ws.Range("D18").Value = "2.212.63"
$ws.Range("E18").Value = "  +2.48%  "
$ws.Range("D19").Value = "41.332.15"
$ws.Range("E19").Value = "  +4.09%  "
$ws.Range("D20").Value = "72.95"
$ws.Range("E20").Value = "  +1.63%  "
$ws.Range("D21").Value = "0.0₃0895"
$ws.Range("E21").Value = "  +5.67%  "
$ws.Range("D22").Value = "6.06"
$ws.Range("E22").Value = "  +0.64%  "
$ws.Range("D23").Value = "251.14"
$ws.Range("E23").Value = "  +10.12%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("E25").Value = "  +1.25%  "
$ws.Range("D26").Value = "2.27"
$ws.Range("E26").Value = "  -2.42%  "
$ws.Range("D27").Value = "9.49"
$ws.Range("E27").Value = "  +0.74%  "
$ws.Range("D28").Value = "168.27"
$ws.Range("E28").Value = "  -2.13%  "
$ws.Range("D29").Value = "0.140"
$ws.Range("E29").Value = "  +0.53%  "
$ws.Range("D30").Value = "19.96"
$ws.Range("E30").Value = "  +1.75%  "
$ws.Range("E31").Value = "  -0.74%  "
$ws.Range("E32").Value = "  -2.05%  "
$ws.Range("E33").Value = "  +0.75%  "
$ws.Range("E34").Value = "  +6.07%  "
$ws.Range("E35").Value = "  +0.93%  "
$ws.Range("E36").Value = "  +1.62%  "
$ws.Range("D37").Value = "6.59"
$ws.Range("E37").Value = "  -5.09%  "
$ws.Range("E38").Value = "  -1.34%  "
$ws.Range("E39").Value = "  -1.64%  "
$ws.Range("D40").Value = "0.000245"
$ws.Range("E40").Value = "  +28.58%  "
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "0.0239"
$ws.Range("E42").Value = "  +5.33%  "
$ws.Range("B43").Value = "FTXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D43").Value = "4.83"
$ws.Range("E43").Value = "  -2.16%  "
$ws.Range("D44").Value = "8.59"
$ws.Range("E44").Value = "  +8.31%  "
$ws.Range("E45").Value = "  +5.63%  "
$ws.Range("D46").Value = "99.02"
$ws.Range("E46").Value = "  -3.57%  "
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("D48").Value = "1.464.93"
$ws.Range("E48").Value = "  -3.16%  "
$ws.Range("E49").Value = "  -5.47%  "
$ws.Range("E50").Value = "  -0.97%  "
$ws.Range("E51").Value = "  -1.26%  "
